$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Agenda table (rows 8-14) ---
# First, copy formatting for cells whose style needs to shift, using source
# cells that already carry the destination style (captured before they are
# overwritten further below).

# F14 needs style 26, currently held by F13
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)

# F11 needs style 21, currently held by F10
$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial(-4122)

# F12 needs style 23, currently held by F11 (but F11 above already changed
# its style) -- use F26 instead, which also carries style 23 and is untouched
$ws.Range("F26").Copy()
$ws.Range("F12").PasteSpecial(-4122)

# F13 needs style 21, use F10 (still style 21) as source
$ws.Range("F10").Copy()
$ws.Range("F13").PasteSpecial(-4122)

# F9 needs to lose its explicit style (reverts to default/general)
$ws.Range("F9").ClearFormats()

$excel.CutCopyMode = $false

# --- Header block (row 3) ---
$ws.Range("E3").Value = " AGENDA Project: Barroc-IT voortgangsvergadering van Projectgroep 6"

# --- Footer notes (rows 17-18) ---
$ws.Range("H17").Value = "Agenda aan :P. Hoek, K. Ly, M. Havermans, F.van Krimpen "
$ws.Range("H18").Value = "Notulen aan :F. van Krimpen, Github, P. Hoek "

# --- Agenda table text/values ---
$ws.Range("F12").Value = "Plan van aanpak bespreken voor inlevering"

# --- Header block (rows 4-5) ---
$ws.Range("B4").Value = "Tijd: 10:00 t/m 10:20"
$ws.Range("B5").Value = "Datum:18-09-2014 "

$ws.Range("F9").Value = "Mededelingen"

# --- Mededelingen note (row 24) ---
$ws.Range("F24").Value = "Goed teamsverband"

# Remaining agenda table edits (rows 8-14)
$ws.Range("C8").Value = 1

$ws.Range("B9").Value = "Voorzitter"
$ws.Range("C9").Value = 2

$ws.Range("B10").Value = "Voorzitter"
$ws.Range("C10").Value = 2
$ws.Range("F10").Value = "Communicatie bij afwezigheid"

$ws.Range("C11").Value = 3
$ws.Range("F11").Value = "Taakverdeling( wie doet wat)"

$ws.Range("C13").Value = 2
$ws.Range("F13").Value = "Vaststellen volgende vergadering"

$ws.Range("B14").Value = "Groep"
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = "t"
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = "Rondvraag"

# --- Selection state ---
[void]$ws.Range("F25").Select()
